# Runtime update: refresh the "as_of_utc" snapshot timestamp on both data
# sheets, and on the "Линейные" (linesmen) sheet add a newly-appeared
# official "Gribovskiy Nikita" (alphabetically sorted into row 10), which
# pushes the existing entries down one row and drops the official who had
# been in the last row ("Sysuev Aleksandr").

$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-10-29 13:09:21"

# --- Sheet 2: "Главные" (main referees) ----------------------------------
# Only the as_of_utc timestamp (column AA) changes, for every data row.
$wsMain = $wb.Worksheets.Item(2)
for ($r = 2; $r -le 26; $r++) {
    $wsMain.Cells.Item($r, 27).Value = $newTimestamp
}

# --- Sheet 3: "Линейные" (linesmen) ---------------------------------------
$wsLine = $wb.Worksheets.Item(3)

# Insert a new row at position 10, shifting the existing rows 10-26 down to
# 11-27 (this preserves every other official's data/types untouched).
$wsLine.Rows.Item(10).Insert()

# The official who used to occupy the last row has now been shifted out of
# the table's range (to row 27); remove that now-duplicate trailing row so
# the sheet keeps its original 26-row extent.
$wsLine.Rows.Item(27).Delete()

# Populate the newly inserted row 10 with the new official's stats.
# Columns: A=Official B=Official_ru C=Games_KHL D=PIM_total E=PIM_home
# F=PIM_away G=PIM_per_game H=PIM_per_game_home I=PIM_per_game_away
# J=PIM_2min_home K=PIM_2min_away L=PIM_5min_home M=PIM_5min_away
# N=PIM_10min_home O=PIM_10min_away P=PIM_20min_home Q=PIM_20min_away
# R=PIM_25min_home S=PIM_25min_away T=PIM_SHB_home U=PIM_SHB_away
# V=PIM_K_home W=PIM_K_away X=PIM_V_home Y=PIM_V_away Z=season_id
$wsLine.Cells.Item(10, 1).Value = "Gribovskiy Nikita"
$wsLine.Cells.Item(10, 2).Value = "Грибовский Никита"
$wsLine.Cells.Item(10, 3).Value = 1
$wsLine.Cells.Item(10, 4).Value = 10
$wsLine.Cells.Item(10, 5).Value = 2
$wsLine.Cells.Item(10, 6).Value = 8
$wsLine.Cells.Item(10, 7).Value = 10
$wsLine.Cells.Item(10, 8).Value = 2
$wsLine.Cells.Item(10, 9).Value = 8
$wsLine.Cells.Item(10, 10).Value = 1
$wsLine.Cells.Item(10, 11).Value = 4
$wsLine.Cells.Item(10, 12).Value = 0
$wsLine.Cells.Item(10, 13).Value = 0
$wsLine.Cells.Item(10, 14).Value = 0
$wsLine.Cells.Item(10, 15).Value = 0
$wsLine.Cells.Item(10, 16).Value = 0
$wsLine.Cells.Item(10, 17).Value = 0
$wsLine.Cells.Item(10, 18).Value = 0
$wsLine.Cells.Item(10, 19).Value = 0
$wsLine.Cells.Item(10, 20).Value = 0
$wsLine.Cells.Item(10, 21).Value = 0
$wsLine.Cells.Item(10, 22).Value = 0
$wsLine.Cells.Item(10, 23).Value = 0
$wsLine.Cells.Item(10, 24).Value = 0
$wsLine.Cells.Item(10, 25).Value = 2
$wsLine.Cells.Item(10, 26).Value = "17"
$wsLine.Cells.Item(10, 27).Value = $newTimestamp

# Refresh the as_of_utc timestamp (column AA) for every other data row.
for ($r = 2; $r -le 26; $r++) {
    if ($r -ne 10) {
        $wsLine.Cells.Item($r, 27).Value = $newTimestamp
    }
}
